# Commit message: "changing FALSE to False"
#
# Column I (rows 2-41) held a boolean formula `=FALSE()` rendered through a
# custom "TRUE"/"FALSE" number format. The edit replaces that with the
# literal text string "False" stored as plain text (number format changed
# to the built-in Text format, @ / numFmtId 49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("I2:I41")

# Switch the column to the Text number format before writing the literal
# word "False" into it.
$rng.NumberFormat = "@"

# Write a formula that evaluates to the text "False" in every cell, then
# convert the whole range to static values so the stored cell becomes a
# plain (non-formula) shared string "False" - exactly like a user typing
# the word into a Text-formatted cell.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 9).Formula = '="False"'
}
$rng.Copy()
$rng.PasteSpecial(-4163)   # xlPasteValues

# A new row was added below the table (row 42) with an empty cell in column
# I carrying the same (Text) formatting as the rest of the column, which is
# also where the selection ended up.
$ws.Cells.Item(41, 9).Copy()
$ws.Cells.Item(42, 9).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(42, 9).ClearContents()
$ws.Range("I42").Select()
